# Daily attendance processing - 2025-11-27 06:37:09
# Reorders the comma-separated "Recorded By" names in column G so that
# "System" is listed first (or moved ahead of the trailing duplicate
# lowercase "system" entry), matching the exact text substitutions:
#   "<name>, System"              -> "System, <name>"
#   "<name>, system, System"      -> "<name>, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    if ($val -eq "backup@backdoor.com, system, System") {
        $cell.Value2 = "backup@backdoor.com, System, system"
    }
    elseif ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "admin@admin.com, System") {
        $cell.Value2 = "System, admin@admin.com"
    }
}
